# "Check forgot password e-mail"
#
# Adds two new reference columns (F, G) to the "ForgotPasswordEmail" sheet
# holding the e-mail body/subject text used by the "forgot password" flow,
# and leaves that sheet as the active tab/selection (mirroring a user who
# just typed these values in and was looking at the new G2 cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ForgotPasswordEmail")

# Make this the active sheet (workbook activeTab moves from IncorrectUserID
# to ForgotPasswordEmail; the old active sheet automatically loses its
# tabSelected flag).
$ws.Activate()

# New header row (row 1) + sample value row (row 2) for columns F and G.
# (Written column-by-column, F then G, so new shared-string entries land in
# the same order as the author's original edit.)
$ws.Range("F1").Value = "Forgot Password E-Mail Body Text"
$ws.Range("F2").Value = "XXX"
$ws.Range("G1").Value = "Forgot Password E-Mail Subject Text"
$ws.Range("G2").Value = "MicroEdge - Automated"

# Match the column sizing of the rest of the sheet (best-fit style widths).
$ws.Columns.Item(6).ColumnWidth = 30.59
$ws.Columns.Item(7).ColumnWidth = 32.74

# Scroll so column D is at the left edge and G2 (the last value typed) is
# the active/selected cell, matching the author's final view.
$ws.Range("G2").Select()
$excel.ActiveWindow.ScrollColumn = 4
